$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 228, shifting existing rows 228-323 down to 229-324
$ws.Rows(228).Insert()

# Populate the new row 228 with the new data record
$ws.Range("A228").Value = 6
$ws.Range("B228").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C228").Value = "Metropolitana"
$ws.Range("D228").Value = 44523
$ws.Range("E228").Value = 13
$ws.Range("F228").Value = 100112039
$ws.Range("G228").Value = "Ciboulette"
$ws.Range("H228").Value = "Sin especificar"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 910
$ws.Range("K228").Value = 700
$ws.Range("L228").Value = 800
$ws.Range("M228").Value = 748
$ws.Range("N228").Value = "`$/docena de atados"
$ws.Range("O228").Value = "Región Metropolitana"
$ws.Range("P228").Value = 249
$ws.Range("Q228").Value = 3
$ws.Range("R228").Value = "Hortaliza"
